$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-typed number formatting on the Price (D) column cells we are
# about to rewrite, so strings like "314.33" are stored as text (matching
# the original inlineStr cells) instead of being auto-coerced to numbers.
$dPriceCells = @("D2","D3","D5","D7","D8","D9","D10","D12","D13","D14","D15","D16","D17","D18","D19","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D48","D51")
foreach ($addr in $dPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Row-by-row value updates (Price / Volume(1h)) ---
$ws.Range("D2").Value = "28.393.38"
$ws.Range("E2").Value = "  +3.94%  "
$ws.Range("D3").Value = "1.795.86"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "314.33"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "0.5447"
$ws.Range("E7").Value = "  +5.86%  "
$ws.Range("D8").Value = "0.3825"
$ws.Range("E8").Value = "  +3.73%  "
$ws.Range("D9").Value = "0.07566"
$ws.Range("E9").Value = "  +2.33%  "
$ws.Range("D10").Value = "42.71"
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("E11").Value = "  +3.30%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").Value = "21.11"
$ws.Range("E13").Value = "  +3.02%  "
$ws.Range("D14").Value = "6.182"
$ws.Range("E14").Value = "  +1.87%  "
$ws.Range("D15").Value = "7.393"
$ws.Range("E15").Value = "  +6.26%  "
$ws.Range("D16").Value = "1.800.40"
$ws.Range("E16").Value = "  +1.91%  "
$ws.Range("D17").Value = "91.47"
$ws.Range("E17").Value = "  +2.49%  "
$ws.Range("D18").Value = "0.00001070"
$ws.Range("E18").Value = "  +2.22%  "
$ws.Range("D19").Value = "0.06448"
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("E21").Value = "  +3.27%  "
$ws.Range("D22").Value = "5.953"
$ws.Range("E22").Value = "  +2.17%  "
$ws.Range("D23").Value = "28.412.53"
$ws.Range("E23").Value = "  +3.92%  "
$ws.Range("D24").Value = "11.45"
$ws.Range("E24").Value = "  +1.84%  "
$ws.Range("D25").Value = "2.129"
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("D26").Value = "159.60"
$ws.Range("E26").Value = "  +3.45%  "
$ws.Range("D27").Value = "20.69"
$ws.Range("E27").Value = "  +2.46%  "
$ws.Range("D28").Value = "2.393"
$ws.Range("E28").Value = "  +2.77%  "
$ws.Range("D29").Value = "2.003.17"
$ws.Range("E29").Value = "  +1.53%  "
$ws.Range("D30").Value = "123.42"
$ws.Range("E30").Value = "  +1.84%  "
$ws.Range("D31").Value = "1.125"
$ws.Range("E31").Value = "  +5.50%  "
$ws.Range("D32").Value = "0.1024"
$ws.Range("E32").Value = "  +4.96%  "
$ws.Range("D33").Value = "5.751"
$ws.Range("E33").Value = "  +3.10%  "
$ws.Range("D34").Value = "3.679"
$ws.Range("E34").Value = "  +1.43%  "
$ws.Range("D37").Value = "0.02322"
$ws.Range("E37").Value = "  +3.45%  "
$ws.Range("D38").Value = "5.156"
$ws.Range("E38").Value = "  +6.56%  "
$ws.Range("D39").Value = "8.766"
$ws.Range("E39").Value = "  +8.04%  "
$ws.Range("D40").Value = "11.66"
$ws.Range("E40").Value = "  +3.65%  "
$ws.Range("D41").Value = "0.6389"
$ws.Range("E41").Value = "  +3.93%  "
$ws.Range("D44").Value = "1.155"
$ws.Range("E44").Value = "  +1.86%  "
$ws.Range("D45").Value = "13.56"
$ws.Range("E45").Value = "  +3.69%  "
$ws.Range("D46").Value = "0.5973"
$ws.Range("E46").Value = "  +3.55%  "
$ws.Range("E47").Value = "  +0.92%  "
$ws.Range("D48").Value = "126.41"
$ws.Range("E48").Value = "  +4.18%  "
$ws.Range("E49").Value = "  +6.18%  "
$ws.Range("E50").Value = "  +3.26%  "
$ws.Range("D51").Value = "0.06961"

# --- Rows 35/36 swapped content (Hedera <-> Algorand) ---
$ws.Range("B35").Value = "Algorand"
$ws.Range("C35").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D35").Value = "0.2353"
$ws.Range("E35").Value = "  +16.45%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.06749"
$ws.Range("E36").Value = "  +13.02%  "

# --- Rows 42/43 swapped content (Frax <-> WEMIXTOKEN) ---
$ws.Range("B42").Value = "WEMIXTOKEN"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "1.417"
$ws.Range("E42").Value = "  -1.22%  "

$ws.Range("B43").Value = "Frax"
$ws.Range("C43").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D43").Value = "0.9999"
$ws.Range("E43").Value = "  -0.09%  "
